# BOT; UPDATE DATA
# Appends one more day of data (2020-05-18 / serial 43969) to the daily
# consultation-count table on sheet "相談件数".
#
# The table previously ended with a footer/note row at row 114
# (shared string index 4). We insert a new row at 114 - which pushes
# that footer row down to row 115 and copies the number formats from
# the row above (row 113) into the freshly inserted row - and then
# fill in the new day's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

# Push the existing footer row (old row 114) down to row 115 and create
# a new, blank row 114 that inherits formatting from row 113.
$ws.Rows(114).Insert()

# Fill in the new day's data in row 114.
$ws.Range("A114").Value = 43969
$ws.Range("B114").Value = 287
$ws.Range("C114").Value = 37965
$ws.Range("D114").Value = 58
$ws.Range("E114").Value = 7642

# Update the active selection to reflect where the user ended up after
# the edit (bottom-right frozen pane, cell B115 - the footer note).
$ws.Activate()
$ws.Range("B115").Select()
